# Update "handback-status.xlsx" timestamps to reflect the generated report times.
$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Latest HO Xliff Generate Date for 5d2f643e-...md
$wsOverview.Range("G2").Value = "2016-08-18 05:04:32"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Correspond Handoff Datetime for 5d2f643e-...zh-cn.xlf
$wsZhCn.Range("H2").Value = "2016-08-18 05:04:27"
# Correspond Handback DateTime for 5d2f643e-...zh-cn.xlf
$wsZhCn.Range("K2").Value = "2016-08-18 05:04:43"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Correspond Handoff Datetime for 5d2f643e-...de-de.xlf (shares value with Overview!G2)
$wsDeDe.Range("H2").Value = "2016-08-18 05:04:32"
# Correspond Handback DateTime for 5d2f643e-...de-de.xlf
$wsDeDe.Range("K2").Value = "2016-08-18 05:04:50"
